$d = $word.ActiveDocument

# 1. Update the date/weather text: replace the trailing "晴123245" with "多云"
$d.Content.Find.Execute("2023年3月26日星期天天气晴123245", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2023年3月26日星期天天气多云", 2)

# 2. Remove the final (empty) paragraph entirely, merging its paragraph mark
#    away so the document body ends with the date paragraph immediately
#    followed by the section properties.
$paras = $d.Paragraphs
$count = $paras.Count
$lastPara = $paras.Item($count)
$secondLastPara = $paras.Item($count - 1)
$r = $d.Range($secondLastPara.Range.End - 1, $lastPara.Range.End)
$r.Delete()
